$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("commands to remember")

$ws2.Range("A4").Value = "mvn dtest"
$ws2.Range("B4").Value = " mvn test -Dtest=AuthenticationControllerUnitTest"
$ws2.Range("C4").Value = "it runs the particular test class , here the class is AuthenticationControllerUnitTest"

$ws2.Range("A4:C4").WrapText = $true

$ws2.Range("C4").Select()

$ws2.Activate()

$wb.Save()
